# "Generate Report for Archive"
# - Update the localization status text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview summary sheet plus
#   each per-locale detail sheet).
# - Narrow the "Latest HO Xliff Generate Date" / per-locale "Status"
#   columns to their new, tighter report width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Overview sheet: columns E (zh-cn) and F (de-de) hold the status value.
if ($wsOverview.Range("E2").Text -eq $oldStatus) {
    $wsOverview.Range("E2").Value = $newStatus
}
if ($wsOverview.Range("F2").Text -eq $oldStatus) {
    $wsOverview.Range("F2").Value = $newStatus
}

# Per-locale detail sheets: column C ("Status") holds the status value.
if ($wsZhCn.Range("C2").Text -eq $oldStatus) {
    $wsZhCn.Range("C2").Value = $newStatus
}
if ($wsDeDe.Range("C2").Text -eq $oldStatus) {
    $wsDeDe.Range("C2").Value = $newStatus
}

# New, narrower column width shared by the resized columns (target stored
# OOXML column width ~= 13.41 characters). ColumnWidth is expressed in
# "characters of the Normal style" and gets snapped to the host's pixel
# grid, so we use the nearest value that lands on the intended width.
$newWidth = 12.5

# Overview: narrow the zh-cn (E) and de-de (F) status columns.
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# zh-cn / de-de detail sheets: narrow the Status column (C).
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
